$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# "Rescatables" (makeup-exam candidates) sheet: grows from 2 data rows to 5.
# Two new students (rows 2-3) plus one more (row 4) are inserted ahead of the
# two students that were already listed, which now move down to rows 5-6.
#
# Columns: A=NC, B=Paterno, C=Materno, D=Nombres, E=Nombre_Largo, F=Grupo, G=Reprobadas
$colA = @(19330051920011, 18330051920281, 18330051920309, 18330051920306, 18330051920429)
$colB = @("DE JESUS", "LOPEZ", "VALENCIA", "TEMOXTLE", "GUERRA")
$colC = @("AGUILAR", "ROSAS", "TORRES", "LARA", "OLMEDO")
$colD = @("HAZIEL", "ALEXA", "LUIS ENRIQUE", "MADAI", "PAOLA BETSABET")
$colE = @("ECOLOGÍA", "TEMAS DE BIOLOGÍA CONTEMPORÁNEA", "TEMAS DE BIOLOGÍA CONTEMPORÁNEA", "TEMAS DE BIOLOGÍA CONTEMPORÁNEA", "TEMAS DE BIOLOGÍA CONTEMPORÁNEA")
$colF = @("4AEM", "6ALCM", "6ALCM", "6ALCM", "6BLCM")
$colG = @(2, 2, 2, 1, 1)

$columns = @($colA, $colB, $colC, $colD, $colE, $colF, $colG)

# Write column by column (all of A, then all of B, ...) so that brand-new
# text values are registered/interned in the same left-to-right, top-to-
# bottom order the source data was entered in.
for ($c = 0; $c -lt $columns.Length; $c++) {
    $colValues = $columns[$c]
    for ($i = 0; $i -lt $colValues.Length; $i++) {
        $r = $i + 2
        $ws.Cells.Item($r, $c + 1).Value2 = $colValues[$i]
    }
}
